# Corrections following third round of review
# The "Materials" sheet has a redundant "subgenus" column (column AS) that
# needs to be removed entirely, shifting every subsequent column one to the
# left (e.g. the former EZ2 becomes EY2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Delete the entire "subgenus" column (column AS).
$ws.Range("AS1").EntireColumn.Delete()
